# Update Products.xlsx with new product data and configurations.
# Adds a new "Avatar" column (H) with a header and avatar file names for
# the first three products.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell, formatted like the other header cells (E1/F1/G1):
# centered horizontally and vertically.
$ws.Range("H1").Value = "Avatar"
$ws.Range("H1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("H1").VerticalAlignment = -4108     # xlCenter

# New avatar values for the first three products, formatted like the rest
# of the data columns (E2:G16): horizontally centered only.
$ws.Range("H2").Value = "girl1"
$ws.Range("H3").Value = "boy1"
$ws.Range("H4").Value = "boy2"
$ws.Range("H2:H4").HorizontalAlignment = -4108  # xlCenter

# Mirror the interactive selection of the whole new column that produced
# this edit.
$ws.Columns.Item(8).Select()
